$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 13
$ws.Range("H13").Value = 881
$ws.Range("I13").Value = 50
$ws.Range("J13").Value = 1296.5
$ws.Range("K13").Value = 50
$ws.Range("L13").Value = 1296.5
$ws.Range("M13").Value = 119
$ws.Range("N13").Value = -1634.5
# Row 34
$ws.Range("H34").Value = 2950
$ws.Range("I34").Value = 2950
$ws.Range("K34").Value = 2950
$ws.Range("M34").Value = -2747
# Row 36
$ws.Range("H36").Value = 2950
$ws.Range("I36").Value = 2950
$ws.Range("K36").Value = 2950
$ws.Range("M36").Value = -2235
# Row 64
$ws.Range("H64").Value = 5000
$ws.Range("I64").Value = 5000
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 5000
$ws.Range("L64").Value = 5000
$ws.Range("M64").Value = -4752
$ws.Range("N64").Value = -5496
# Row 67
$ws.Range("H67").Value = 5000
$ws.Range("I67").Value = 5000
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 5000
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = -4142
$ws.Range("N67").Value = -6716
# Row 70
$ws.Range("H70").Value = 54999.5
$ws.Range("I70").Value = 10000
$ws.Range("K70").Value = 30000
$ws.Range("M70").Value = -29730
# Row 73
$ws.Range("H73").Value = 54999.5
$ws.Range("I73").Value = 10000
$ws.Range("K73").Value = 30000
$ws.Range("M73").Value = -29064
# Row 138
$ws.Range("H138").Value = 3258.3225
$ws.Range("J138").Value = 3406.3635
$ws.Range("L138").Value = 10219.0905
$ws.Range("N138").Value = -20499.0905

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value = 13344.667
$ws.Range("I37").Value = 34
$ws.Range("K37").Value = 34
$ws.Range("M37").Value = 239
# Row 55
$ws.Range("H55").Value = 19500
# Row 93
$ws.Range("H93").Value = 18500
$ws.Range("I93").Value = 12000
$ws.Range("J93").Value = 25000
$ws.Range("K93").Value = 12000
$ws.Range("L93").Value = 25000
$ws.Range("M93").Value = -9504
$ws.Range("N93").Value = -29992
# Row 102
$ws.Range("H102").Value = 3000
$ws.Range("J102").Value = 3000
$ws.Range("L102").Value = 3000
$ws.Range("N102").Value = -6244

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 167
$ws.Range("I5").Value = 167
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 167
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -54
$ws.Range("N5").Value = $null
# Row 7
$ws.Range("H7").Value = 6667000
$ws.Range("I7").Value = 6667000
$ws.Range("K7").Value = 6667000
$ws.Range("M7").Value = -6666887
# Row 12
$ws.Range("H12").Value = 4184.6665
$ws.Range("I12").Value = 5002
$ws.Range("J12").Value = 3776
$ws.Range("K12").Value = 5002
$ws.Range("L12").Value = 3776
$ws.Range("M12").Value = -4834
$ws.Range("N12").Value = -4112
# Row 19
$ws.Range("H19").Value = 105
$ws.Range("I19").Value = 105
$ws.Range("K19").Value = 105
$ws.Range("M19").Value = 68
# Row 86
$ws.Range("H86").Value = 1123.75
$ws.Range("I86").Value = 1069.25
$ws.Range("J86").Value = 1287.25
$ws.Range("K86").Value = 1069.25
$ws.Range("L86").Value = 1287.25
$ws.Range("M86").Value = 53.75
$ws.Range("N86").Value = -3533.25
# Row 89
$ws.Range("H89").Value = 1123.75
$ws.Range("I89").Value = 1069.25
$ws.Range("J89").Value = 1287.25
$ws.Range("K89").Value = 5346.25
$ws.Range("L89").Value = 6436.25
$ws.Range("M89").Value = 269.75
$ws.Range("N89").Value = -17668.25
# Row 92
$ws.Range("H92").Value = 32499.5
$ws.Range("J92").Value = 32499.5
$ws.Range("L92").Value = 32499.5
$ws.Range("N92").Value = -37491.5
# Row 95
$ws.Range("H95").Value = 10499.833
$ws.Range("J95").Value = 10499.833
$ws.Range("L95").Value = 10499.833
$ws.Range("N95").Value = -15991.833

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 104
$ws.Range("I2").Value = 104
$ws.Range("K2").Value = 104
$ws.Range("M2").Value = 9
# Row 21
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 1000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 1000
$ws.Range("M21").Value = $null
$ws.Range("N21").Value = -1470
# Row 41
$ws.Range("H41").Value = 20000
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").Value = $null
# Row 60
$ws.Range("H60").Value = 21850
$ws.Range("I60").Value = 7250
$ws.Range("J60").Value = 25500
$ws.Range("K60").Value = 7250
$ws.Range("L60").Value = 25500
$ws.Range("M60").Value = -6739
$ws.Range("N60").Value = -26522
# Row 62
$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
# Row 65
$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240
# Row 132
$ws.Range("H132").Value = 750
$ws.Range("I132").Value = 750
$ws.Range("K132").Value = 2250
$ws.Range("M132").Value = 280

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = $null
# Row 29
$ws.Range("H29").Value = 200
$ws.Range("J29").Value = 200
$ws.Range("L29").Value = 600
$ws.Range("N29").Value = -1154

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = $null
# Row 46
$ws.Range("H46").Value = 11105.125
$ws.Range("I46").Value = 2210.25
$ws.Range("K46").Value = 2210.25
$ws.Range("M46").Value = -2054.25
# Row 54
$ws.Range("H54").Value = 10000
$ws.Range("J54").Value = 10000
$ws.Range("L54").Value = 10000
$ws.Range("N54").Value = -10780
# Row 57
$ws.Range("H57").Value = 28750
$ws.Range("J57").Value = 28750
$ws.Range("L57").Value = 28750
$ws.Range("N57").Value = -30390

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 12
$ws.Range("H12").Value = 1832.6666
$ws.Range("I12").Value = 2499
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 2499
$ws.Range("L12").Value = 500
$ws.Range("M12").Value = -2329
$ws.Range("N12").Value = -840
# Row 26
$ws.Range("H26").Value = 10000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 10000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 10000
$ws.Range("M26").Value = $null
$ws.Range("N26").Value = -10590
# Row 40
$ws.Range("H40").Value = 6589
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").Value = $null
# Row 95
$ws.Range("H95").Value = 25000
$ws.Range("J95").Value = 25000
$ws.Range("L95").Value = 25000
$ws.Range("N95").Value = -30492

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 13
$ws.Range("H13").Value = 950
$ws.Range("I13").Value = 900
$ws.Range("K13").Value = 900
$ws.Range("M13").Value = -760
# Row 81
$ws.Range("H81").Value = 9028.286
$ws.Range("I81").Value = 9699.666999999999
$ws.Range("K81").Value = 19399.334
$ws.Range("M81").Value = -18338.334
# Row 84
$ws.Range("H84").Value = 9028.286
$ws.Range("I84").Value = 9699.666999999999
$ws.Range("K84").Value = 96996.67
$ws.Range("M84").Value = -91692.67
# Row 132
$ws.Range("H132").Value = 8604.817999999999
$ws.Range("I132").Value = 8078.1113
$ws.Range("J132").Value = 10975
$ws.Range("K132").Value = 24234.3339
$ws.Range("L132").Value = 32925
$ws.Range("M132").Value = -21704.3339
$ws.Range("N132").Value = -37985
